$wb = $excel.ActiveWorkbook

# Sheet2: update B2 value from 30 to 25
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("B2").Value = 25

# Sheet2: change selection to C1 (no longer the active/tabSelected sheet)
$ws2.Range("C1").Select()

# Sheet1: becomes the active sheet / selected tab
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()

$wb.Save()
